$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- Add header row 10 (copy format from row 4, then set the same header text) ---
$ws.Range("A4:J4").Copy() | Out-Null
$ws.Range("A10").PasteSpecial(-4122) | Out-Null

$ws.Range("A10").Value = "TC_ID"
$ws.Range("B10").Value = "TestcaseName"
$ws.Range("C10").Value = "OpportunityName"
$ws.Range("D10").Value = "RelatedTo"
$ws.Range("E10").Value = "ContactName"
$ws.Range("F10").Value = "Subject"
$ws.Range("G10").Value = "OrganizationName"
$ws.Range("H10").Value = "BillingAddress"
$ws.Range("I10").Value = "ProductName"
$ws.Range("J10").Value = "Qty"

# --- Add data row 11 (copy format from row 5, then set the new test-case data) ---
$ws.Range("A5:J5").Copy() | Out-Null
$ws.Range("A11").PasteSpecial(-4122) | Out-Null

$ws.Range("A11").Value = "TC_016"
$ws.Range("B11").Value = "Add_Quote_To_Opportunity"
$ws.Range("C11").Value = "Client2"
$ws.Range("D11").Value = "Contacts"
$ws.Range("E11").Value = "Asha89"

# --- Update J5: turn the numeric 1 into a quote-prefixed text "1" ---
$ws.Range("J5").Value = "'1"

$ws.Range("F11").Value = "Check"
$ws.Range("G11").Value = "Instagram_31"
$ws.Range("H11").Value = "3rd floor Gopalan coworks, `nKathriguppe"
$ws.Range("I11").Value = "Volvo"
$ws.Range("J11").Value = "'1"

$ws.Range("F11").Borders.Item(7).LineStyle = -4142

$ws.Rows.Item(11).RowHeight = 46.5

# --- Column F: widen and stop relying on bestFit ---
$ws.Columns.Item(6).ColumnWidth = 25.43

# --- Selection matches the author's last recorded position ---
$ws.Range("G15").Select() | Out-Null
